$wb = $excel.ActiveWorkbook

# "SolverSettings" is the sheet that gets the new RPS / capacity-credit input row.
$ws = $wb.Worksheets.Item("SolverSettings")

# Add the new scenario-flag row (row 10): include_RPS, defaulted to "N" for every scenario column B:G.
$ws.Range("A10").Value = "include_RPS"
$ws.Range("B10:G10").Value = "N"

# Make SolverSettings the active sheet/tab and leave the selection on H10,
# matching where the cursor ends up after entering the new row of data.
$ws.Activate()
$ws.Range("H10").Select()
